# Insert a new column before column DL ("Note") to hold the new
# "DemonstrationProjectIdentifier" field. This shifts every existing
# column from DL onward one position to the right (DL->DM, DM->DN, ...,
# MN->MO), which also grows the sheet's used range from MN5 to MO5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("DL:DL").Insert()

# Populate the header for the newly inserted column.
$ws.Range("DL1").Value2 = "DemonstrationProjectIdentifier"

# Update the record identifier that changed for this example row set.
$ws.Range("A2").Value2 = "690148897e79911955eafb87"
$ws.Range("A3").Value2 = "690148897e79911955eafb87"
$ws.Range("A4").Value2 = "690148897e79911955eafb87"
$ws.Range("A5").Value2 = "690148897e79911955eafb87"
